$d = $word.ActiveDocument

# The document has a single section whose header/footer "logo" pictures
# need their OOXML `name` swapped:
#   - the BTec logo pictures (currently named "image2.jpg") become "image1.jpg"
#   - the Pearson logo pictures (currently named "image1.png") become "image2.png"
# (the picture content / relationship itself is untouched - only the
# cosmetic drawing name changes).

$sec = $d.Sections.Item(1)

# Headers: primary (Item 1) and first-page (Item 2) headers both carry the
# BTec logo inline picture.
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -ge 1) {
        $ish = $hdr.Range.InlineShapes.Item(1)
        if ($ish.AlternativeText -eq "BTec_Logo-Orange") {
            $ish.Name = "image1.jpg"
        }
    }
}

# Footers: primary (Item 1) and first-page (Item 2) footers both carry the
# Pearson logo inline picture.
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -ge 1) {
        $ish = $ftr.Range.InlineShapes.Item(1)
        if ($ish.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $ish.Name = "image2.png"
        }
    }
}
